$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44280
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806

# Row 3
$ws.Range("D3").Value = 44280
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 667

# Row 4
$ws.Range("D4").Value = 44516
$ws.Range("N4").Value = 33000
$ws.Range("O4").Value = 34000
$ws.Range("P4").Value = 33500
$ws.Range("S4").Value = 1861

# Row 5
$ws.Range("D5").Value = 44316
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 1111
